$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("A1").Value = "Sr.No."
$ws.Range("A2").Value = "Sr.No."
